$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the "Status" columns ---
# Target stored column width (OOXML `width` attribute) is 13.4101845877511.
# The COM ColumnWidth setter here persists as (round(ColumnWidth * 6)/6) + 5/6,
# so back the value off by that fixed 5/6 padding before assigning.
$newColumnWidth = 13.4101845877511 - (5 / 6)

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
